$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1879.8334
$ws.Range("I2").Value = 1879.8334
$ws.Range("K2").Value = 1879.8334
$ws.Range("M2").Value = -1766.8334
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = $null
$ws.Range("H19").Value = 3113.647
$ws.Range("I19").Value = 180.14285
$ws.Range("J19").Value = 5167.1
$ws.Range("K19").Value = 180.14285
$ws.Range("L19").Value = 5167.1
$ws.Range("M19").Value = -5.14285000000001
$ws.Range("N19").Value = -5517.1
$ws.Range("H33").Value = 807.5
$ws.Range("I33").Value = 602.2857
$ws.Range("J33").Value = 1525.75
$ws.Range("K33").Value = 602.2857
$ws.Range("L33").Value = 1525.75
$ws.Range("M33").Value = -373.2857
$ws.Range("N33").Value = -1983.75
$ws.Range("H51").Value = 2285
$ws.Range("I51").Value = 2072.7273
$ws.Range("J51").Value = 2544.4443
$ws.Range("K51").Value = 2072.7273
$ws.Range("L51").Value = 2544.4443
$ws.Range("M51").Value = -1588.7273
$ws.Range("N51").Value = -3512.4443
$ws.Range("H137").Value = 2129931
$ws.Range("I137").Value = 3450724
$ws.Range("J137").Value = 1986.6666
$ws.Range("K137").Value = 10352172
$ws.Range("L137").Value = 5959.9998
$ws.Range("M137").Value = -10349622
$ws.Range("N137").Value = -11059.9998

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 800
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("H61").Value = 1845.9412
$ws.Range("I61").Value = 815.67346
$ws.Range("J61").Value = 4502.9473
$ws.Range("K61").Value = 815.67346
$ws.Range("L61").Value = 4502.9473
$ws.Range("M61").Value = -603.67346
$ws.Range("N61").Value = -4926.9473
$ws.Range("H132").Value = 17547570
$ws.Range("I132").Value = 27781882
$ws.Range("J132").Value = 3038
$ws.Range("K132").Value = 83345646
$ws.Range("L132").Value = 9114
$ws.Range("M132").Value = -83343116
$ws.Range("N132").Value = -14174
$ws.Range("H136").Value = 1845.9412
$ws.Range("I136").Value = 815.67346
$ws.Range("J136").Value = 4502.9473
$ws.Range("K136").Value = 2447.02038
$ws.Range("L136").Value = 13508.8419
$ws.Range("M136").Value = 102.9796200000001
$ws.Range("N136").Value = -18608.8419

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 800
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H134").Value = 3305.0466
$ws.Range("I134").Value = 3042.484
$ws.Range("J134").Value = 3983.3333
$ws.Range("K134").Value = 9127.451999999999
$ws.Range("L134").Value = 11949.9999
$ws.Range("M134").Value = -6592.451999999999
$ws.Range("N134").Value = -17019.9999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2155.8572
$ws.Range("I22").Value = 272.75
$ws.Range("K22").Value = 272.75
$ws.Range("M22").Value = 77.25
$ws.Range("H31").Value = 1353446.2
$ws.Range("I31").Value = 2382223
$ws.Range("J31").Value = 3176.8438
$ws.Range("K31").Value = 2382223
$ws.Range("L31").Value = 3176.8438
$ws.Range("M31").Value = -2381928
$ws.Range("N31").Value = -3766.8438
$ws.Range("H34").Value = 1353446.2
$ws.Range("I34").Value = 2382223
$ws.Range("J34").Value = 3176.8438
$ws.Range("K34").Value = 2382223
$ws.Range("L34").Value = 3176.8438
$ws.Range("M34").Value = -2382021
$ws.Range("N34").Value = -3580.8438
$ws.Range("H132").Value = 2985.5
$ws.Range("I132").Value = 1903.5769
$ws.Range("J132").Value = 4743.625
$ws.Range("K132").Value = 5710.7307
$ws.Range("L132").Value = 14230.875
$ws.Range("M132").Value = -3180.7307
$ws.Range("N132").Value = -19290.875
$ws.Range("H134").Value = 1556.5283
$ws.Range("I134").Value = 1269.92
$ws.Range("J134").Value = 6333.3335
$ws.Range("K134").Value = 3809.76
$ws.Range("L134").Value = 19000.0005
$ws.Range("M134").Value = -1274.76
$ws.Range("N134").Value = -24070.0005

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6540
$ws.Range("J34").Value = 7728.75
$ws.Range("L34").Value = 23186.25
$ws.Range("N34").Value = -23354.25
$ws.Range("H39").Value = 2836.4167
$ws.Range("J39").Value = 3253.7
$ws.Range("L39").Value = 9761.099999999999
$ws.Range("N39").Value = -10349.1
$ws.Range("H55").Value = 3300
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 1500
$ws.Range("M55").Value = -1323
$ws.Range("H68").Value = 1883.9857
$ws.Range("I68").Value = 623.4
$ws.Range("J68").Value = 3144.5715
$ws.Range("K68").Value = 1870.2
$ws.Range("L68").Value = 9433.7145
$ws.Range("M68").Value = -1059.2
$ws.Range("N68").Value = -11055.7145
$ws.Range("H71").Value = 1883.9857
$ws.Range("I71").Value = 623.4
$ws.Range("J71").Value = 3144.5715
$ws.Range("K71").Value = 5610.599999999999
$ws.Range("L71").Value = 28301.1435
$ws.Range("M71").Value = -1554.599999999999
$ws.Range("N71").Value = -36413.1435
$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -16118
$ws.Range("H106").Value = 3871
$ws.Range("J106").Value = 3871
$ws.Range("L106").Value = 11613
$ws.Range("N106").Value = -13505

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2071.2917
$ws.Range("I102").Value = 1428
$ws.Range("K102").Value = 1428
$ws.Range("M102").Value = 194
$ws.Range("H132").Value = 3762.7646
$ws.Range("I132").Value = 2679.1365
$ws.Range("J132").Value = 5749.4165
$ws.Range("K132").Value = 8037.4095
$ws.Range("L132").Value = 17248.2495
$ws.Range("M132").Value = -5507.4095
$ws.Range("N132").Value = -22308.2495
$ws.Range("H135").Value = 29000
$ws.Range("J135").Value = 29000
$ws.Range("L135").Value = 29000
$ws.Range("N135").Value = -39140

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 4995
$ws.Range("I18").Value = 4995
$ws.Range("K18").Value = 4995
$ws.Range("M18").Value = -4823
$ws.Range("H20").Value = 14251
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 14251
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 14251
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -14703
$ws.Range("H132").Value = 2244.7354
$ws.Range("I132").Value = 1220.72
$ws.Range("J132").Value = 5089.222
$ws.Range("K132").Value = 3662.16
$ws.Range("L132").Value = 15267.666
$ws.Range("M132").Value = -1132.16
$ws.Range("N132").Value = -20327.666
$ws.Range("H136").Value = 2085620.6
$ws.Range("I136").Value = 2859050
$ws.Range("J136").Value = 3310.7693
$ws.Range("K136").Value = 8577150
$ws.Range("L136").Value = 9932.3079
$ws.Range("M136").Value = -8574600
$ws.Range("N136").Value = -15032.3079

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1384.1333
$ws.Range("I136").Value = 627.8461
$ws.Range("J136").Value = 6300
$ws.Range("K136").Value = 1883.5383
$ws.Range("L136").Value = 18900
$ws.Range("M136").Value = 666.4617000000001
$ws.Range("N136").Value = -24000
